$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1229.3684
$ws.Range("I19").Value = 1186
$ws.Range("J19").Value = 1254.6666
$ws.Range("K19").Value = 1186
$ws.Range("L19").Value = 1254.6666
$ws.Range("M19").Value = -1011
$ws.Range("N19").Value = -1604.6666

$ws.Range("H64").Value = 4190.8184
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 4262.375
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 4262.375
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -4758.375

$ws.Range("H67").Value = 4190.8184
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 4262.375
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 4262.375
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -5978.375

$ws.Range("H107").Value = 1260
$ws.Range("I107").Value = 1705
$ws.Range("J107").Value = 370
$ws.Range("K107").Value = 1705
$ws.Range("L107").Value = 370
$ws.Range("M107").Value = 215
$ws.Range("N107").Value = -4210

$ws.Range("H116").Value = 596473.5600000001
$ws.Range("I116").Value = 12940
$ws.Range("J116").Value = 839612.5600000001
$ws.Range("K116").Value = 12940
$ws.Range("L116").Value = 839612.5600000001
$ws.Range("M116").Value = -9498
$ws.Range("N116").Value = -846496.5600000001

$ws.Range("H129").Value = 1114.3043
$ws.Range("J129").Value = 1469.8572
$ws.Range("L129").Value = 4409.571599999999
$ws.Range("N129").Value = -14409.5716

$ws.Range("H137").Value = 332383.72
$ws.Range("I137").Value = 568424.3
$ws.Range("J137").Value = 1926.9333
$ws.Range("K137").Value = 1705272.9
$ws.Range("L137").Value = 5780.7999
$ws.Range("M137").Value = -1702722.9
$ws.Range("N137").Value = -10880.7999

$ws.Range("H138").Value = 3151.134
$ws.Range("I138").Value = 2137.5715
$ws.Range("J138").Value = 3431.1973
$ws.Range("K138").Value = 6412.7145
$ws.Range("L138").Value = 10293.5919
$ws.Range("M138").Value = -1272.7145
$ws.Range("N138").Value = -20573.5919

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 886.1818
$ws.Range("I74").Value = 791.7646999999999
$ws.Range("J74").Value = 986.5
$ws.Range("K74").Value = 791.7646999999999
$ws.Range("L74").Value = 986.5
$ws.Range("M74").Value = 82.23530000000005
$ws.Range("N74").Value = -2734.5

$ws.Range("H77").Value = 886.1818
$ws.Range("I77").Value = 791.7646999999999
$ws.Range("J77").Value = 986.5
$ws.Range("K77").Value = 3958.8235
$ws.Range("L77").Value = 4932.5
$ws.Range("M77").Value = 409.1765000000005
$ws.Range("N77").Value = -13668.5

$ws.Range("H97").Value = 317.77777
$ws.Range("I97").Value = 317.77777
$ws.Range("K97").Value = 317.77777
$ws.Range("M97").Value = 178.22223

$ws.Range("H122").Value = 55556920
$ws.Range("I122").Value = 90910070
$ws.Range("J122").Value = 1957
$ws.Range("K122").Value = 272730210
$ws.Range("L122").Value = 5871
$ws.Range("M122").Value = -272727760
$ws.Range("N122").Value = -10771

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1049
$ws.Range("I99").Value = 1054.4445
$ws.Range("K99").Value = 1054.4445
$ws.Range("M99").Value = 443.5554999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15605.412
$ws.Range("I31").Value = 1595.8572
$ws.Range("J31").Value = 25412.1
$ws.Range("K31").Value = 1595.8572
$ws.Range("L31").Value = 25412.1
$ws.Range("M31").Value = -1300.8572
$ws.Range("N31").Value = -26002.1

$ws.Range("H34").Value = 15605.412
$ws.Range("I34").Value = 1595.8572
$ws.Range("J34").Value = 25412.1
$ws.Range("K34").Value = 1595.8572
$ws.Range("L34").Value = 25412.1
$ws.Range("M34").Value = -1393.8572
$ws.Range("N34").Value = -25816.1

$ws.Range("H132").Value = 6806387.5
$ws.Range("I132").Value = 11112274
$ws.Range("J132").Value = 7619.4736
$ws.Range("K132").Value = 33336822
$ws.Range("L132").Value = 22858.4208
$ws.Range("M132").Value = -33334292
$ws.Range("N132").Value = -27918.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 789.5897
$ws.Range("I122").Value = 763.3939
$ws.Range("J122").Value = 933.6667
$ws.Range("K122").Value = 6870.5451
$ws.Range("L122").Value = 8403.0003
$ws.Range("M122").Value = -4420.5451
$ws.Range("N122").Value = -13303.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 23259408
$ws.Range("I132").Value = 32259716
$ws.Range("J132").Value = 8611.666999999999
$ws.Range("K132").Value = 96779148
$ws.Range("L132").Value = 25835.001
$ws.Range("M132").Value = -96776618
$ws.Range("N132").Value = -30895.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7400.4
$ws.Range("I68").Value = 10999.667
$ws.Range("J68").Value = 2001.5
$ws.Range("K68").Value = 10999.667
$ws.Range("L68").Value = 2001.5
$ws.Range("M68").Value = -10250.667
$ws.Range("N68").Value = -3499.5

$ws.Range("H71").Value = 7400.4
$ws.Range("I71").Value = 10999.667
$ws.Range("J71").Value = 2001.5
$ws.Range("K71").Value = 54998.335
$ws.Range("L71").Value = 10007.5
$ws.Range("M71").Value = -51254.335
$ws.Range("N71").Value = -17495.5

$ws.Range("H87").Value = 100000000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0

$ws.Range("H90").Value = 100000000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0

$ws.Range("H122").Value = 312501500
$ws.Range("I122").Value = 333335330
$ws.Range("J122").Value = 250000000
$ws.Range("K122").Value = 1000005990
$ws.Range("L122").Value = 750000000
$ws.Range("M122").Value = -1000003540
$ws.Range("N122").Value = -750004900

$ws.Range("H136").Value = 3443.8704
$ws.Range("I136").Value = 4531.8823
$ws.Range("K136").Value = 13595.6469
$ws.Range("M136").Value = -11045.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5557755.5
$ws.Range("I81").Value = 11111633
$ws.Range("J81").Value = 3877.7778
$ws.Range("K81").Value = 22223266
$ws.Range("L81").Value = 7755.5556
$ws.Range("M81").Value = -22222205
$ws.Range("N81").Value = -9877.5556

$ws.Range("H84").Value = 5557755.5
$ws.Range("I84").Value = 11111633
$ws.Range("J84").Value = 3877.7778
$ws.Range("K84").Value = 111116330
$ws.Range("L84").Value = 38777.778
$ws.Range("M84").Value = -111111026
$ws.Range("N84").Value = -49385.778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()